# Update "Name of Algo" column (header "F", spreadsheet column E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 12.987
$ws.Range("E14").Value = 13.239
$ws.Range("E21").Value = 13.449
$ws.Range("E23").Value = 13.136
$ws.Range("E25").Value = 12.659
